{"js": "// Add a bookmark named \"CourseContent2\" around the existing \"CourseContent2\"\n// text in the Course Module table (the author had bookmarked every other\n// placeholder in that row except this one; this change adds the missing\n// bookmark). Inserting the bookmark automatically shifts the w:id values of\n// every bookmark that follows it later in the document, matching the\n// renumbering seen in the diff.\nconst searchResults = context.document.body.search(\"CourseContent2\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find \"CourseContent2\" text in the document.');\n}\n\nsearchResults.items[0].insertBookmark(\"CourseContent2\");\nawait context.sync();\n", "ps1": "# Add a bookmark named \"CourseContent2\" around the existing \"CourseContent2\"\n# text in the Course Module table (the author had bookmarked every other\n# placeholder in that row except this one; this change adds the missing\n# bookmark). Inserting the bookmark automatically shifts the w:id values of\n# every bookmark that follows it later in the document, matching the\n# renumbering seen in the diff.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"CourseContent2\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find \"CourseContent2\" text in the document.'\n}\n\n$d.Bookmarks.Add(\"CourseContent2\", $range)\n"}
